$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"

$ws.Range('D2').Value = '26.062.15'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '1.646.32'
$ws.Range('E3').Value = '  -1.36%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = '215.61'
$ws.Range('E5').Value = '  +2.35%  '
$ws.Range('D6').Value = '0.5216'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('D8').Value = '0.2609'
$ws.Range('E8').Value = '  -0.71%  '
$ws.Range('D9').Value = '0.06363'
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('D10').Value = '20.83'
$ws.Range('E10').Value = '  -1.60%  '
$ws.Range('D11').Value = '0.07673'
$ws.Range('E11').Value = '  +1.78%  '
$ws.Range('D12').Value = '1.654.06'
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('E13').Value = '  -0.27%  '
$ws.Range('D14').Value = '1.868.00'
$ws.Range('E14').Value = '  -1.61%  '
$ws.Range('D15').Value = '0.5546'
$ws.Range('E15').Value = '  +0.97%  '
$ws.Range('D16').Value = '0.0₅8262'
$ws.Range('E16').Value = '  +3.10%  '
$ws.Range('D17').Value = '64.97'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('D18').Value = '26.080.17'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('E20').Value = '  -0.30%  '
$ws.Range('D21').Value = '188.61'
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('E22').Value = '  -1.03%  '
$ws.Range('D23').Value = '6.239'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').Value = '146.20'
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('D26').Value = '0.1218'
$ws.Range('E26').Value = '  -1.72%  '
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').Value = '15.83'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').Value = '1.384'
$ws.Range('E29').Value = '  +2.51%  '
$ws.Range('D30').Value = '0.05950'
$ws.Range('E30').Value = '  -5.60%  '
$ws.Range('D31').Value = '1.269'
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('D32').Value = '3.403'
$ws.Range('E32').Value = '  -0.37%  '
$ws.Range('D33').Value = '3.400'
$ws.Range('E33').Value = '  -3.38%  '
$ws.Range('D34').Value = '1.663'
$ws.Range('E34').Value = '  +0.98%  '
$ws.Range('D35').Value = '0.9975'
$ws.Range('E35').Value = '  -0.77%  '
$ws.Range('D36').Value = '2.392'
$ws.Range('E36').Value = '  -0.39%  '
$ws.Range('D37').Value = '2.754'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('D38').Value = '0.5618'
$ws.Range('E38').Value = '  -6.82%  '
$ws.Range('D40').Value = '5.849'
$ws.Range('E40').Value = '  -3.97%  '
$ws.Range('D41').Value = '0.8556'
$ws.Range('E41').Value = '  -0.91%  '
$ws.Range('D43').Value = '1.030.28'
$ws.Range('E43').Value = '  -7.65%  '
$ws.Range('D44').Value = '99.29'
$ws.Range('E44').Value = '  -1.20%  '
$ws.Range('D45').Value = '1.795.59'
$ws.Range('E45').Value = '  -1.44%  '
$ws.Range('D46').Value = '0.0₈112'
$ws.Range('E46').Value = '  +5.06%  '
$ws.Range('D47').Value = '55.83'
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('D48').Value = '1.003'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').Value = '8.095'
$ws.Range('E49').Value = '  +0.08%  '
$ws.Range('E50').Value = '  -1.63%  '
$ws.Range('E51').Value = '  -0.49%  '

$ws.Range('D5').Style = "Normal"
$ws.Range('D6').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D28').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D34').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D49').Style = "Normal"
